$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 33
$ws.Range("F3").Value = 116
$ws.Range("G3").Value = 242
$ws.Range("H3").Value = 214
$ws.Range("I3").Value = 261
$ws.Range("J3").Value = 188
$ws.Range("K3").Value = 172
$ws.Range("L3").Value = 29
$ws.Range("E4").Value = 240
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 62
$ws.Range("H4").Value = 115
$ws.Range("I4").Value = 23
$ws.Range("J4").Value = 46
$ws.Range("K4").Value = 99
$ws.Range("L4").Value = 239
$ws.Range("E5").Value = 16602
$ws.Range("F5").Value = 79904
$ws.Range("G5").Value = 404651
$ws.Range("H5").Value = 580478
$ws.Range("I5").Value = 389079
$ws.Range("J5").Value = 754174
$ws.Range("K5").Value = 3365151
$ws.Range("L5").Value = 15334
$ws.Range("E6").Value = 1748671
$ws.Range("F6").Value = 954
$ws.Range("G6").Value = 128148
$ws.Range("H6").Value = 2735201
$ws.Range("I6").Value = 22029
$ws.Range("J6").Value = 4118116
$ws.Range("K6").Value = 104670
$ws.Range("L6").Value = 195624
$ws.Range("E7").Value = 467
$ws.Range("F7").Value = 464
$ws.Range("H7").Value = 463
$ws.Range("I7").Value = 461
$ws.Range("J7").Value = 461
$ws.Range("K7").Value = 466
$ws.Range("L7").Value = 469
$ws.Range("E8").Value = 462
$ws.Range("F8").Value = 477
$ws.Range("G8").Value = 463
$ws.Range("H8").Value = 467
$ws.Range("I8").Value = 469
$ws.Range("J8").Value = 466
$ws.Range("K8").Value = 463
$ws.Range("L8").Value = 462
$ws.Range("E9").Value = 818
$ws.Range("F9").Value = 6734
$ws.Range("H9").Value = 26689
$ws.Range("I9").Value = 16527
$ws.Range("J9").Value = 408297
$ws.Range("L9").Value = 1859
$ws.Range("E10").Value = 907406
$ws.Range("F10").Value = 477
$ws.Range("G10").Value = 96784
$ws.Range("I10").Value = 9063
$ws.Range("J10").Value = 3808470
$ws.Range("K10").Value = 14510
$ws.Range("L10").Value = 36436
$ws.Range("E11").Value = 503.0909090909091
$ws.Range("F11").Value = 688.8275862068965
$ws.Range("G11").Value = 1672.111570247934
$ws.Range("H11").Value = 2712.514018691589
$ws.Range("I11").Value = 1490.724137931034
$ws.Range("J11").Value = 4011.563829787234
$ws.Range("K11").Value = 19564.83139534884
$ws.Range("L11").Value = 528.7586206896551
$ws.Range("E12").Value = 7286.129166666667
$ws.Range("F12").Value = 477
$ws.Range("G12").Value = 2066.903225806452
$ws.Range("H12").Value = 23784.35652173913
$ws.Range("I12").Value = 957.7826086956521
$ws.Range("J12").Value = 89524.26086956522
$ws.Range("K12").Value = 1057.272727272727
$ws.Range("L12").Value = 818.510460251046
$ws.Range("E13").Value = 470
$ws.Range("F13").Value = 476
$ws.Range("G13").Value = 469.25
$ws.Range("H13").Value = 471
$ws.Range("I13").Value = 476
$ws.Range("E14").Value = 470
$ws.Range("F14").Value = 477
$ws.Range("H14").Value = 473
$ws.Range("I14").Value = 470
$ws.Range("E15").Value = 479
$ws.Range("F15").Value = 477
$ws.Range("G15").Value = 476
$ws.Range("H15").Value = 481
$ws.Range("I15").Value = 480
$ws.Range("J15").Value = 471
$ws.Range("K15").Value = 479
$ws.Range("E16").Value = 470
$ws.Range("F16").Value = 477
$ws.Range("G16").Value = 476
$ws.Range("H16").Value = 587
$ws.Range("I16").Value = 477
$ws.Range("K16").Value = 477
$ws.Range("L16").Value = 471
$ws.Range("E17").Value = 488
$ws.Range("F17").Value = 478
$ws.Range("G17").Value = 588.75
$ws.Range("H17").Value = 2203.5
$ws.Range("I17").Value = 1194
$ws.Range("J17").Value = 1120.5
$ws.Range("K17").Value = 648.5
$ws.Range("L17").Value = 478
$ws.Range("E18").Value = 479
$ws.Range("F18").Value = 477
$ws.Range("G18").Value = 480
$ws.Range("H18").Value = 2986
$ws.Range("I18").Value = 479
$ws.Range("K18").Value = 746.5
$ws.Range("L18").Value = 480
$ws.Range("E19").Value = 18
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 119.5
$ws.Range("H19").Value = 1732.5
$ws.Range("I19").Value = 718
$ws.Range("J19").Value = 650.5
$ws.Range("K19").Value = 178.5
$ws.Range("L19").Value = 8
$ws.Range("E20").Value = 9
$ws.Range("F20").ClearContents()
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 2513
$ws.Range("I20").Value = 9
$ws.Range("K20").Value = 276.5
$ws.Range("L20").Value = 10
$ws.Range("E25").Value = 32
$ws.Range("F25").Value = 115
$ws.Range("G25").Value = 241
$ws.Range("H25").Value = 213
$ws.Range("I25").Value = 260
$ws.Range("J25").Value = 187
$ws.Range("K25").Value = 171
$ws.Range("L25").Value = 28
$ws.Range("E26").Value = 239
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 61
$ws.Range("H26").Value = 114
$ws.Range("I26").Value = 22
$ws.Range("J26").Value = 45
$ws.Range("K26").Value = 98
$ws.Range("L26").Value = 238
$ws.Range("E27").Value = 134019
$ws.Range("F27").Value = 80276
$ws.Range("G27").Value = 225491
$ws.Range("H27").Value = 229176
$ws.Range("I27").Value = 230792
$ws.Range("J27").Value = 200653
$ws.Range("K27").Value = 210588
$ws.Range("L27").Value = 197355
$ws.Range("E28").Value = 144248
$ws.Range("F28").Value = 233
$ws.Range("G28").Value = 198562
$ws.Range("H28").Value = 223077
$ws.Range("I28").Value = 219744
$ws.Range("J28").Value = 192445
$ws.Range("K28").Value = 213955
$ws.Range("L28").Value = 212586
$ws.Range("E29").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("G29").Value = -57
$ws.Range("H29").ClearContents()
$ws.Range("I29").ClearContents()
$ws.Range("J29").ClearContents()
$ws.Range("K29").ClearContents()
$ws.Range("L29").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("F30").Value = 233
$ws.Range("G30").ClearContents()
$ws.Range("H30").ClearContents()
$ws.Range("I30").ClearContents()
$ws.Range("J30").ClearContents()
$ws.Range("K30").ClearContents()
$ws.Range("L30").ClearContents()
$ws.Range("E31").Value = 35637
$ws.Range("F31").Value = 20160
$ws.Range("G31").Value = 16106
$ws.Range("H31").Value = 30070
$ws.Range("I31").Value = 31410
$ws.Range("J31").Value = 31910
$ws.Range("K31").Value = 41718
$ws.Range("L31").Value = 109783
$ws.Range("E32").Value = 21955
$ws.Range("F32").Value = 233
$ws.Range("G32").Value = 42468
$ws.Range("H32").Value = 38333
$ws.Range("I32").Value = 72522
$ws.Range("J32").Value = 38870
$ws.Range("K32").Value = 51579
$ws.Range("L32").Value = 46455
$ws.Range("E33").Value = 4188.09375
$ws.Range("F33").Value = 698.0521739130435
$ws.Range("G33").Value = 935.6473029045643
$ws.Range("H33").Value = 1075.943661971831
$ws.Range("I33").Value = 887.6615384615385
$ws.Range("J33").Value = 1073.010695187166
$ws.Range("K33").Value = 1231.508771929825
$ws.Range("L33").Value = 7048.392857142857
$ws.Range("E34").Value = 603.5481171548117
$ws.Range("F34").Value = 233
$ws.Range("G34").Value = 3255.114754098361
$ws.Range("H34").Value = 1956.815789473684
$ws.Range("I34").Value = 9988.363636363636
$ws.Range("J34").Value = 4276.555555555556
$ws.Range("K34").Value = 2183.214285714286
$ws.Range("L34").Value = 893.2184873949579
$ws.Range("E35").Value = 1.75
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 1
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = 4.5
$ws.Range("K35").Value = 2
$ws.Range("L35").Value = 1
$ws.Range("E36").Value = 3
$ws.Range("F36").Value = 233
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 1.25
$ws.Range("I36").Value = 4
$ws.Range("J36").Value = 13
$ws.Range("K36").Value = 1
$ws.Range("L36").Value = 1
$ws.Range("E37").Value = 6
$ws.Range("F37").Value = 2
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 5
$ws.Range("I37").Value = 2
$ws.Range("J37").Value = 57
$ws.Range("K37").Value = 4
$ws.Range("L37").Value = 6
$ws.Range("E38").Value = 46
$ws.Range("F38").Value = 233
$ws.Range("G38").Value = 36
$ws.Range("H38").Value = 6
$ws.Range("I38").Value = 386
$ws.Range("J38").Value = 265
$ws.Range("K38").Value = 4
$ws.Range("L38").Value = 2
$ws.Range("E39").Value = 3672.25
$ws.Range("F39").Value = 7
$ws.Range("G39").Value = 213
$ws.Range("H39").Value = 438
$ws.Range("I39").Value = 8
$ws.Range("J39").Value = 676.5
$ws.Range("K39").Value = 66.5
$ws.Range("L39").Value = 1657.5
$ws.Range("E40").Value = 174.5
$ws.Range("F40").Value = 233
$ws.Range("G40").Value = 1300
$ws.Range("H40").Value = 636.5
$ws.Range("I40").Value = 11211.75
$ws.Range("J40").Value = 2840
$ws.Range("K40").Value = 205.25
$ws.Range("L40").Value = 28
$ws.Range("E41").Value = 3670.5
$ws.Range("F41").Value = 6
$ws.Range("G41").Value = 212
$ws.Range("H41").Value = 437
$ws.Range("I41").Value = 7
$ws.Range("J41").Value = 672
$ws.Range("K41").Value = 64.5
$ws.Range("L41").Value = 1656.5
$ws.Range("E42").Value = 171.5
$ws.Range("F42").ClearContents()
$ws.Range("G42").Value = 1299
$ws.Range("H42").Value = 635.25
$ws.Range("I42").Value = 11207.75
$ws.Range("J42").Value = 2827
$ws.Range("K42").Value = 204.25
$ws.Range("L42").Value = 27
